# "Waiting to add screenshots." — merge the two runs of the third bullet
# on the Conclusions slide ("Learned about deploying " + "web applications")
# into a single run, and drop the stray trailing endParaRPr mark.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Rebuild the text frame with only the first two (unchanged) bullets; this
# drops the third paragraph entirely, taking its lingering endParaRPr with it.
$tr.Text = "Learned how to develop full scale web application`rLearned how to represent data with JSON"

# Re-append the third bullet as a brand-new paragraph/run, so it is written
# out as a single <a:r> (no leftover endParaRPr) with the correct run format.
[void]$tr.InsertAfter("`rLearned about deploying web applications")
